# Processed Results - nexus 5x until https request
# Replace the gyroscope run data in column B (B2:B31) with the newly
# captured "Nexus 5X" measurements, which ripples through every dependent
# formula (E3/E4 average, D7/E7 min/max, D10/E10 quartiles, D13/E13 IQR,
# D16/E16 stdev / relative stdev) automatically via recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$newValues = @(
    132.108192,
    134.26883999999899,
    118.526975999999,
    114.360012,
    114.20568,
    114.823008,
    114.668675999999,
    114.97734,
    115.13167199999999,
    114.360012,
    114.668675999999,
    115.90333200000001,
    114.668675999999,
    117.909648,
    114.360012,
    114.668675999999,
    114.360012,
    113.89701599999999,
    115.13167199999999,
    114.823008,
    114.668675999999,
    115.286003999999,
    115.440336,
    114.360012,
    114.051348,
    114.668675999999,
    113.742683999999,
    113.89701599999999,
    114.51434399999999,
    114.20568
)

$row = 2
foreach ($v in $newValues) {
    $ws.Cells.Item($row, 2).Value = $v
    $row++
}

# The two hidden helper defined names that back the embedded charts get a
# second generation (v1.2 / v1.3) pointing at the same ranges as v1.0/v1.1.
$wb.Names.Add("_xlchart.v1.2", $ws.Range("A2:A31"))
$wb.Names.Add("_xlchart.v1.3", $ws.Range("B2:B31"))
$wb.Names.Item("_xlchart.v1.2").Visible = $false
$wb.Names.Item("_xlchart.v1.3").Visible = $false

# The value axis of the scatter chart no longer pins a fixed maximum of
# 105 - let it scale automatically with the refreshed data.
$chart = $ws.ChartObjects().Item(1).Chart
$valueAxis = $chart.Axes().Item(2)
$valueAxis.MaximumScaleIsAuto = $true

# Selection moved from the old "last touched" cell (E16) to the refreshed
# data range, with B4 as the active cell.
$ws.Activate()
$ws.Range("B4:B31,B2,B3").Select()
